$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Update the status text on every worksheet that carries it (Overview has the
# per-language status duplicated in columns E/F, zh-cn and de-de each carry
# it once in column C).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find($oldStatus)
    $guard = 0
    while ($found -ne $null -and $guard -lt 100) {
        $found.Value = $newStatus
        $found = $used.Find($oldStatus)
        $guard = $guard + 1
    }
}

# The status column(s) were narrowed after the shorter text was written
# ("Ready for handoff" -> "In Translation"). Resize them accordingly.
$newStatusColWidth = 12.58

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E1").ColumnWidth = $newStatusColWidth
$ws1.Range("F1").ColumnWidth = $newStatusColWidth

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C1").ColumnWidth = $newStatusColWidth

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C1").ColumnWidth = $newStatusColWidth
